$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6454346776008606
$ws.Range("B1").Value = 1.767122507095337
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.935716390609741
$ws.Range("E1").Value = 1.137291669845581
